$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 2024 class-date / assignment-date refresh ---------------------------

# Week 1 row
$ws.Range("B2").Value = "Jan. 8"
$ws.Range("C2").Value = "A Common Vocabulary                                    "
$ws.Range("E2").Value = "Student survey (Jan. 10)"

# Week 2 row
$ws.Range("B3").Value = "Online"

# Week 3 row
$ws.Range("B4").Value = "Jan. 22"
$ws.Range("E4").Value = "DARE #1 (Jan. 21)"

# Week 4 row
$ws.Range("B5").Value = "Jan. 29"
$ws.Range("D5").Value = "MM Ch. 9 <br> [Angrist & Lavy 1999](https://doi.org/10.1162/003355399556061) <br> [Dee & Penner 2017](https://journals.sagepub.com/doi/full/10.3102/0002831216677002)"

# Week 5 row
$ws.Range("B6").Value = "Feb. 5"
$ws.Range("E6").Value = "DARE #2 (Feb. 5)"

# Week 6 row
$ws.Range("B7").Value = "Feb. 12"

# Week 7 row
$ws.Range("B8").Value = "Feb. 19"
$ws.Range("E8").Value = "DARE #3 (Feb. 18)"

# Week 8 row
$ws.Range("B9").Value = "Feb. 26"

# Week 9 row
$ws.Range("B10").Value = "Mar. 4"
$ws.Range("E10").Value = "DARE #4 (Mar. 3)"

# Week 10 row
$ws.Range("B11").Value = "Mar. 11"
$ws.Range("E11").Value = "Research project presentation (Mar. 11)"

# Finals row
$ws.Range("B12").Value = "NA"
$ws.Range("E12").Value = "Final research project (Mar. 20)"

# --- formatting cleanup (drop now-unused "applyAlignment" cell formats) --

$ws.Range("E2:E11").NumberFormat = "@"
$ws.Range("E12").ClearFormats()

# --- restore selection ----------------------------------------------------

$ws.Range("C2").Select()
